$d = $word.ActiveDocument

# --- Edit 1: paragraph 1 - insert the "Content Moved to Integrated Doc" text box
#     as a new leading run (AlternateContent drawing) before the existing run.
$para1 = $d.Paragraphs(1)
$rng1 = $para1.Range
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="en-US"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="57797D83" wp14:editId="6AD0F55C"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>0</wp:posOffset></wp:positionV><wp:extent cx="1828800" cy="1828800"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="5" name="Text Box 5"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1828800" cy="1828800"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln><a:effectLst/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="262626" w:themeColor="text1" w:themeTint="D9"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:shadow w14:blurRad="0" w14:dist="38100" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="bl"><w14:schemeClr w14:val="accent5"/></w14:shadow><w14:textOutline w14:w="6731" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="bg1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="262626" w:themeColor="text1" w:themeTint="D9"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:shadow w14:blurRad="0" w14:dist="38100" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="bl"><w14:schemeClr w14:val="accent5"/></w14:shadow><w14:textOutline w14:w="6731" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="bg1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:t>Content Moved to Integrated Doc</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="none" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype w14:anchorId="57797D83" id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="Text Box 5" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;margin-left:0;margin-top:0;width:2in;height:2in;z-index:251659264;visibility:visible;mso-wrap-style:none;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQDXMe3YJAIAAFUEAAAOAAAAZHJzL2Uyb0RvYy54bWysVE2P2jAQvVfqf7B8LwEELRsRVnRXVJXQ&#10;7kpQ7dk4NolkeyzbkNBf37GTsHTbU9WLM18ez8x7k+V9qxU5C+drMAWdjMaUCMOhrM2xoD/2m08L&#10;SnxgpmQKjCjoRXh6v/r4YdnYXEyhAlUKRzCJ8XljC1qFYPMs87wSmvkRWGHQKcFpFlB1x6x0rMHs&#10;WmXT8fhz1oArrQMuvEfrY+ekq5RfSsHDs5ReBKIKirWFdLp0HuKZrZYsPzpmq5r3ZbB/qEKz2uCj&#10;11SPLDBycvUfqXTNHXiQYcRBZyBlzUXqAbuZjN91s6uYFakXHI631zH5/5eWP51fHKnLgs4pMUwj&#10;RHvRBvIVWjKP02mszzFoZzEstGhGlAe7R2NsupVOxy+2Q9CPc75cZxuT8XhpMV0sxuji6BsUzJ+9&#10;XbfOh28CNIlCQR2Cl2bKzlsfutAhJL5mYFMrlQBU5jcD5uwsIjGgvx076SqOUmgPbd/eAcoLdueg&#10;Y4e3fFNjBVvmwwtzSAesGikenvGQCpqCQi9RUoH7+Td7jEeU0EtJg/QqqEH+U6K+G0TvbjKbRTYm&#10;ZTb/MkXF3XoOtx5z0g+A/J3gKlmexBgf1CBKB/oV92Ad30QXMxxfLmgYxIfQUR73iIv1OgUh/ywL&#10;W7OzPKaOA4zT3bevzNkegoDoPcFAQ5a/Q6KLjTe9XZ8C4pFgiuPtZorwRgW5m4Du9ywux62eot7+&#10;BqtfAAAA//8DAFBLAwQUAAYACAAAACEAS4kmzdYAAAAFAQAADwAAAGRycy9kb3ducmV2LnhtbEyP&#10;0U7DMAxF35H4h8hIvLF0FaBSmk5owDMw+ACvMU1p41RNthW+HoOQxovlq2tdn1utZj+oPU2xC2xg&#10;uchAETfBdtwaeHt9vChAxYRscQhMBj4pwqo+PamwtOHAL7TfpFZJCMcSDbiUxlLr2DjyGBdhJBbv&#10;PUwek8ip1XbCg4T7QedZdq09diwfHI60dtT0m503UGT+qe9v8ufoL7+WV259Hx7GD2POz+a7W1CJ&#10;5nQ8hh98QYdamLZhxzaqwYAUSb9TvLwoRG7/Fl1X+j99/Q0AAP//AwBQSwECLQAUAAYACAAAACEA&#10;toM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQA&#10;BgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQItABQA&#10;BgAIAAAAIQDXMe3YJAIAAFUEAAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBLAQIt&#10;ABQABgAIAAAAIQBLiSbN1gAAAAUBAAAPAAAAAAAAAAAAAAAAAH4EAABkcnMvZG93bnJldi54bWxQ&#10;SwUGAAAAAAQABADzAAAAgQUAAAAA&#10;" filled="f" stroked="f"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="262626" w:themeColor="text1" w:themeTint="D9"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:shadow w14:blurRad="0" w14:dist="38100" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="bl"><w14:schemeClr w14:val="accent5"/></w14:shadow><w14:textOutline w14:w="6731" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="bg1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="262626" w:themeColor="text1" w:themeTint="D9"/><w:sz w:val="72"/><w:szCs w:val="72"/><w14:shadow w14:blurRad="0" w14:dist="38100" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="bl"><w14:schemeClr w14:val="accent5"/></w14:shadow><w14:textOutline w14:w="6731" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:solidFill><w14:schemeClr w14:val="bg1"/></w14:solidFill><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:t>Content Moved to Integrated Doc</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p></w:txbxContent></v:textbox></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Model-Based Standard &#8211; a standard based on and published as a model that can be reused or implemented directly within other software</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng1.InsertXML($xml1)

# --- Edit 2: paragraph 4 - fix "Continuous Exploratoin" -> split into two runs
#     wrapped with proofErr spell-check markers (text unchanged, still misspelled).
$para4 = $d.Paragraphs(4)
$rng4 = $para4.Range
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Continuous </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Exploratoin</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng4.InsertXML($xml4)
